# link01.xlsx: add two HYPERLINK() formula cells to Sheet1 (E1, E2) that
# exercise workbook-qualified and sheet-qualified external-style references,
# per the Workbook::parse_range / Cell::formula / Cell::get_hyperlink work
# described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New cell E1 (same row as the existing A1/C1 headers): a HYPERLINK formula
# that references another sheet in the (external-looking) workbook.
$ws.Range("E1").Formula = '=HYPERLINK("[link01.xlsx]Sheet1!A5","ByFormula")'

# New cell E2 (row 2, which previously had no populated cells): a HYPERLINK
# formula using a bare (sheet-less) reference.
$ws.Range("E2").Formula = '=HYPERLINK("[link01.xlsx]A5","ByFormula2")'

# Match the formatting already used by the other header-row cells (A1/C1).
$ws.Range("E1").Style = $ws.Range("A1").Style
$ws.Range("E2").Style = $ws.Range("A1").Style

# The new selection/active cell ends up on E2.
$ws.Range("E2").Select()

# Best-effort: the workbook's built-in hyperlink cell style was authored
# under a French locale ("Lien hypertexte"); normalize it to the English
# built-in name ("Hyperlink"). No-op/safe if the host doesn't support
# renaming a registered cell style.
try {
    $style = $wb.Styles.Item("Lien hypertexte")
    $style.Name = "Hyperlink"
} catch {
}
